$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 38 (2024-10-01 data row); row 39 (blank footer) shifts up to become row 38
$ws.Rows.Item(38).Delete()

# Update header info cells
$ws.Range("B4").Value2 = "2024-11-01 ~ 2024-11-30"
$ws.Range("B5").Value2 = "2024년 12월 08일 16시 45분 59초"

# Update the daily visitor data rows (8-37) with November 2024 values
$ws.Range("A8").Value2 = "2024-11-30"
$ws.Range("B8").Value2 = "토"
$ws.Range("C8").Value2 = "9"
$ws.Range("D8").Value2 = "0"
$ws.Range("E8").Value2 = "0"
$ws.Range("F8").Value2 = "9"

$ws.Range("A9").Value2 = "2024-11-29"
$ws.Range("B9").Value2 = "금"
$ws.Range("C9").Value2 = "6"
$ws.Range("D9").Value2 = "0"
$ws.Range("E9").Value2 = "0"
$ws.Range("F9").Value2 = "6"

$ws.Range("A10").Value2 = "2024-11-28"
$ws.Range("B10").Value2 = "목"
$ws.Range("C10").Value2 = "8"
$ws.Range("D10").Value2 = "0"
$ws.Range("E10").Value2 = "0"
$ws.Range("F10").Value2 = "8"

$ws.Range("A11").Value2 = "2024-11-27"
$ws.Range("B11").Value2 = "수"
$ws.Range("C11").Value2 = "8"
$ws.Range("D11").Value2 = "0"
$ws.Range("E11").Value2 = "0"
$ws.Range("F11").Value2 = "8"

$ws.Range("A12").Value2 = "2024-11-26"
$ws.Range("B12").Value2 = "화"
$ws.Range("C12").Value2 = "11"
$ws.Range("D12").Value2 = "0"
$ws.Range("E12").Value2 = "0"
$ws.Range("F12").Value2 = "11"

$ws.Range("A13").Value2 = "2024-11-25"
$ws.Range("B13").Value2 = "월"
$ws.Range("C13").Value2 = "9"
$ws.Range("D13").Value2 = "0"
$ws.Range("E13").Value2 = "0"
$ws.Range("F13").Value2 = "9"

$ws.Range("A14").Value2 = "2024-11-24"
$ws.Range("B14").Value2 = "일"
$ws.Range("C14").Value2 = "5"
$ws.Range("D14").Value2 = "0"
$ws.Range("E14").Value2 = "0"
$ws.Range("F14").Value2 = "5"

$ws.Range("A15").Value2 = "2024-11-23"
$ws.Range("B15").Value2 = "토"
$ws.Range("C15").Value2 = "12"
$ws.Range("D15").Value2 = "0"
$ws.Range("E15").Value2 = "1"
$ws.Range("F15").Value2 = "11"

$ws.Range("A16").Value2 = "2024-11-22"
$ws.Range("B16").Value2 = "금"
$ws.Range("C16").Value2 = "7"
$ws.Range("D16").Value2 = "0"
$ws.Range("E16").Value2 = "0"
$ws.Range("F16").Value2 = "7"

$ws.Range("A17").Value2 = "2024-11-21"
$ws.Range("B17").Value2 = "목"
$ws.Range("C17").Value2 = "7"
$ws.Range("D17").Value2 = "0"
$ws.Range("E17").Value2 = "0"
$ws.Range("F17").Value2 = "7"

$ws.Range("A18").Value2 = "2024-11-20"
$ws.Range("B18").Value2 = "수"
$ws.Range("C18").Value2 = "7"
$ws.Range("D18").Value2 = "1"
$ws.Range("E18").Value2 = "0"
$ws.Range("F18").Value2 = "6"

$ws.Range("A19").Value2 = "2024-11-19"
$ws.Range("B19").Value2 = "화"
$ws.Range("C19").Value2 = "7"
$ws.Range("D19").Value2 = "0"
$ws.Range("E19").Value2 = "0"
$ws.Range("F19").Value2 = "7"

$ws.Range("A20").Value2 = "2024-11-18"
$ws.Range("B20").Value2 = "월"
$ws.Range("C20").Value2 = "7"
$ws.Range("D20").Value2 = "2"
$ws.Range("E20").Value2 = "0"
$ws.Range("F20").Value2 = "5"

$ws.Range("A21").Value2 = "2024-11-17"
$ws.Range("B21").Value2 = "일"
$ws.Range("C21").Value2 = "10"
$ws.Range("D21").Value2 = "0"
$ws.Range("E21").Value2 = "0"
$ws.Range("F21").Value2 = "10"

$ws.Range("A22").Value2 = "2024-11-16"
$ws.Range("B22").Value2 = "토"
$ws.Range("C22").Value2 = "10"
$ws.Range("D22").Value2 = "0"
$ws.Range("E22").Value2 = "0"
$ws.Range("F22").Value2 = "10"

$ws.Range("A23").Value2 = "2024-11-15"
$ws.Range("B23").Value2 = "금"
$ws.Range("C23").Value2 = "14"
$ws.Range("D23").Value2 = "4"
$ws.Range("E23").Value2 = "0"
$ws.Range("F23").Value2 = "10"

$ws.Range("A24").Value2 = "2024-11-14"
$ws.Range("B24").Value2 = "목"
$ws.Range("C24").Value2 = "8"
$ws.Range("D24").Value2 = "3"
$ws.Range("E24").Value2 = "0"
$ws.Range("F24").Value2 = "5"

$ws.Range("A25").Value2 = "2024-11-13"
$ws.Range("B25").Value2 = "수"
$ws.Range("C25").Value2 = "14"
$ws.Range("D25").Value2 = "1"
$ws.Range("E25").Value2 = "0"
$ws.Range("F25").Value2 = "13"

$ws.Range("A26").Value2 = "2024-11-12"
$ws.Range("B26").Value2 = "화"
$ws.Range("C26").Value2 = "6"
$ws.Range("D26").Value2 = "1"
$ws.Range("E26").Value2 = "0"
$ws.Range("F26").Value2 = "5"

$ws.Range("A27").Value2 = "2024-11-11"
$ws.Range("B27").Value2 = "월"
$ws.Range("C27").Value2 = "11"
$ws.Range("D27").Value2 = "0"
$ws.Range("E27").Value2 = "0"
$ws.Range("F27").Value2 = "11"

$ws.Range("A28").Value2 = "2024-11-10"
$ws.Range("B28").Value2 = "일"
$ws.Range("C28").Value2 = "12"
$ws.Range("D28").Value2 = "1"
$ws.Range("E28").Value2 = "1"
$ws.Range("F28").Value2 = "10"

$ws.Range("A29").Value2 = "2024-11-09"
$ws.Range("B29").Value2 = "토"
$ws.Range("C29").Value2 = "12"
$ws.Range("D29").Value2 = "0"
$ws.Range("E29").Value2 = "1"
$ws.Range("F29").Value2 = "11"

$ws.Range("A30").Value2 = "2024-11-08"
$ws.Range("B30").Value2 = "금"
$ws.Range("C30").Value2 = "14"
$ws.Range("D30").Value2 = "4"
$ws.Range("E30").Value2 = "0"
$ws.Range("F30").Value2 = "11"

$ws.Range("A31").Value2 = "2024-11-07"
$ws.Range("B31").Value2 = "목"
$ws.Range("C31").Value2 = "9"
$ws.Range("D31").Value2 = "2"
$ws.Range("E31").Value2 = "0"
$ws.Range("F31").Value2 = "7"

$ws.Range("A32").Value2 = "2024-11-06"
$ws.Range("B32").Value2 = "수"
$ws.Range("C32").Value2 = "5"
$ws.Range("D32").Value2 = "0"
$ws.Range("E32").Value2 = "0"
$ws.Range("F32").Value2 = "5"

$ws.Range("A33").Value2 = "2024-11-05"
$ws.Range("B33").Value2 = "화"
$ws.Range("C33").Value2 = "7"
$ws.Range("D33").Value2 = "0"
$ws.Range("E33").Value2 = "0"
$ws.Range("F33").Value2 = "7"

$ws.Range("A34").Value2 = "2024-11-04"
$ws.Range("B34").Value2 = "월"
$ws.Range("C34").Value2 = "8"
$ws.Range("D34").Value2 = "0"
$ws.Range("E34").Value2 = "0"
$ws.Range("F34").Value2 = "8"

$ws.Range("A35").Value2 = "2024-11-03"
$ws.Range("B35").Value2 = "일"
$ws.Range("C35").Value2 = "7"
$ws.Range("D35").Value2 = "0"
$ws.Range("E35").Value2 = "0"
$ws.Range("F35").Value2 = "7"

$ws.Range("A36").Value2 = "2024-11-02"
$ws.Range("B36").Value2 = "토"
$ws.Range("C36").Value2 = "10"
$ws.Range("D36").Value2 = "0"
$ws.Range("E36").Value2 = "0"
$ws.Range("F36").Value2 = "10"

$ws.Range("A37").Value2 = "2024-11-01"
$ws.Range("B37").Value2 = "금"
$ws.Range("C37").Value2 = "4"
$ws.Range("D37").Value2 = "0"
$ws.Range("E37").Value2 = "0"
$ws.Range("F37").Value2 = "4"

